$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 7.620274999999999
$ws.Range("H2").Value = 22.860825
$ws.Range("I2").Value = 0.6584612850834004
$ws.Range("J2").Value = 0.6584612850834003
$ws.Range("M2").Value = 6.101885666666667
$ws.Range("N2").Value = 18.305657
$ws.Range("O2").Value = 0.1093737608697887
$ws.Range("P2").Value = 0.1093737608697887
$ws.Range("Q2").Value = 46.49804679855833
$ws.Range("R2").Value = 418.482421187025
$ws.Range("S2").Value = 0.07201838713672563
$ws.Range("T2").Value = 0.07201838713672562
# Row 3
$ws.Range("G3").Value = 7.620274999999999
$ws.Range("H3").Value = 22.860825
$ws.Range("I3").Value = 0.6584612850834004
$ws.Range("J3").Value = 0.6584612850834003
$ws.Range("N3").Value = 87.53628900000001
$ws.Range("O3").Value = 0.5230171820937495
$ws.Range("P3").Value = 0.5230171820937495
$ws.Range("Q3").Value = 222.350198219825
$ws.Range("R3").Value = 2001.151783978425
$ws.Range("S3").Value = 0.3443865658421491
$ws.Range("T3").Value = 0.3443865658421491
# Row 4
$ws.Range("G4").Value = 7.620274999999999
$ws.Range("H4").Value = 22.860825
$ws.Range("I4").Value = 0.6584612850834004
$ws.Range("J4").Value = 0.6584612850834003
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.146644
$ws.Range("N4").Value = 0.439932
$ws.Range("O4").Value = 0.002628532664354407
$ws.Range("P4").Value = 0.002628532664354407
$ws.Range("Q4").Value = 1.1174676071
$ws.Range("R4").Value = 10.0572084639
$ws.Range("S4").Value = 0.001730786996054497
$ws.Range("T4").Value = 0.001730786996054497
# Row 5
$ws.Range("G5").Value = 7.620274999999999
$ws.Range("H5").Value = 22.860825
$ws.Range("I5").Value = 0.6584612850834004
$ws.Range("J5").Value = 0.6584612850834003
$ws.Range("M5").Value = 15.02284966666667
$ws.Range("N5").Value = 45.068549
$ws.Range("O5").Value = 0.2692783275177917
$ws.Range("P5").Value = 0.2692783275177917
$ws.Range("Q5").Value = 114.4782457436583
$ws.Range("R5").Value = 1030.304211692925
$ws.Range("S5").Value = 0.1773093535824739
$ws.Range("T5").Value = 0.1773093535824739
# Row 6
$ws.Range("G6").Value = 7.620274999999999
$ws.Range("H6").Value = 22.860825
$ws.Range("I6").Value = 0.6584612850834004
$ws.Range("J6").Value = 0.6584612850834003
$ws.Range("M6").Value = 5.288900666666667
$ws.Range("N6").Value = 15.866702
$ws.Range("O6").Value = 0.09480134312252211
$ws.Range("P6").Value = 0.09480134312252211
$ws.Range("Q6").Value = 40.30287752768333
$ws.Range("R6").Value = 362.72589774915
$ws.Range("S6").Value = 0.06242301422008829
$ws.Range("T6").Value = 0.06242301422008829
# Row 7
$ws.Range("G7").Value = 7.620274999999999
$ws.Range("H7").Value = 22.860825
$ws.Range("I7").Value = 0.6584612850834004
$ws.Range("J7").Value = 0.6584612850834003
$ws.Range("M7").Value = 0.050258
$ws.Range("N7").Value = 0.150774
$ws.Range("O7").Value = 0.0009008537317934847
$ws.Range("P7").Value = 0.0009008537317934848
$ws.Range("Q7").Value = 0.38297978095
$ws.Range("R7").Value = 3.44681802855
$ws.Range("S7").Value = 0.0005931773059089149
$ws.Range("T7").Value = 0.0005931773059089148
# Row 8
$ws.Range("I8").Value = 0.262323813236933
$ws.Range("J8").Value = 0.262323813236933
$ws.Range("M8").Value = 6.101885666666667
$ws.Range("N8").Value = 18.305657
$ws.Range("O8").Value = 0.1093737608697887
$ws.Range("P8").Value = 0.1093737608697887
$ws.Range("Q8").Value = 18.52431603890311
$ws.Range("R8").Value = 166.718844350128
$ws.Range("S8").Value = 0.02869134201942744
$ws.Range("T8").Value = 0.02869134201942744
# Row 9
$ws.Range("I9").Value = 0.262323813236933
$ws.Range("J9").Value = 0.262323813236933
$ws.Range("N9").Value = 87.53628900000001
$ws.Range("O9").Value = 0.5230171820937495
$ws.Range("P9").Value = 0.5230171820937495
$ws.Range("Q9").Value = 88.58190024585068
$ws.Range("R9").Value = 797.2371022126562
$ws.Range("S9").Value = 0.1371998615952677
$ws.Range("T9").Value = 0.1371998615952677
# Row 10
$ws.Range("I10").Value = 0.262323813236933
$ws.Range("J10").Value = 0.262323813236933
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.146644
$ws.Range("N10").Value = 0.439932
$ws.Range("O10").Value = 0.002628532664354407
$ws.Range("P10").Value = 0.002628532664354407
$ws.Range("Q10").Value = 0.4451869388586667
$ws.Range("R10").Value = 4.006682449728
$ws.Range("S10").Value = 0.0006895267117312834
$ws.Range("T10").Value = 0.0006895267117312836
# Row 11
$ws.Range("I11").Value = 0.262323813236933
$ws.Range("J11").Value = 0.262323813236933
$ws.Range("M11").Value = 15.02284966666667
$ws.Range("N11").Value = 45.068549
$ws.Range("O11").Value = 0.2692783275177917
$ws.Range("P11").Value = 0.2692783275177917
$ws.Range("Q11").Value = 45.60688781018844
$ws.Range("R11").Value = 410.461990291696
$ws.Range("S11").Value = 0.07063811769653087
$ws.Range("T11").Value = 0.07063811769653089
# Row 12
$ws.Range("I12").Value = 0.262323813236933
$ws.Range("J12").Value = 0.262323813236933
$ws.Range("M12").Value = 5.288900666666667
$ws.Range("N12").Value = 15.866702
$ws.Range("O12").Value = 0.09480134312252211
$ws.Range("P12").Value = 0.09480134312252211
$ws.Range("Q12").Value = 16.05622799242311
$ws.Range("R12").Value = 144.506051931808
$ws.Range("S12").Value = 0.0248686498278829
$ws.Range("T12").Value = 0.0248686498278829
# Row 13
$ws.Range("I13").Value = 0.262323813236933
$ws.Range("J13").Value = 0.262323813236933
$ws.Range("M13").Value = 0.050258
$ws.Range("N13").Value = 0.150774
$ws.Range("O13").Value = 0.0009008537317934847
$ws.Range("P13").Value = 0.0009008537317934848
$ws.Range("Q13").Value = 0.1525749786773333
$ws.Range("R13").Value = 1.373174808096
$ws.Range("S13").Value = 0.0002363153860927883
$ws.Range("T13").Value = 0.0002363153860927883
# Row 14
$ws.Range("G14").Value = 0.9167423333333334
$ws.Range("H14").Value = 2.750227
$ws.Range("I14").Value = 0.07921490167966665
$ws.Range("J14").Value = 0.07921490167966663
$ws.Range("M14").Value = 6.101885666666667
$ws.Range("N14").Value = 18.305657
$ws.Range("O14").Value = 0.1093737608697887
$ws.Range("P14").Value = 0.1093737608697887
$ws.Range("Q14").Value = 5.593856903793223
$ws.Range("R14").Value = 50.344712134139
$ws.Range("S14").Value = 0.008664031713635686
$ws.Range("T14").Value = 0.008664031713635684
# Row 15
$ws.Range("G15").Value = 0.9167423333333334
$ws.Range("H15").Value = 2.750227
$ws.Range("I15").Value = 0.07921490167966665
$ws.Range("J15").Value = 0.07921490167966663
$ws.Range("N15").Value = 87.53628900000001
$ws.Range("O15").Value = 0.5230171820937495
$ws.Range("P15").Value = 0.5230171820937495
$ws.Range("Q15").Value = 26.74940727640034
$ws.Range("R15").Value = 240.744665487603
$ws.Range("S15").Value = 0.04143075465633267
$ws.Range("T15").Value = 0.04143075465633266
# Row 16
$ws.Range("G16").Value = 0.9167423333333334
$ws.Range("H16").Value = 2.750227
$ws.Range("I16").Value = 0.07921490167966665
$ws.Range("J16").Value = 0.07921490167966663
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.146644
$ws.Range("N16").Value = 0.439932
$ws.Range("O16").Value = 0.002628532664354407
$ws.Range("P16").Value = 0.002628532664354407
$ws.Range("Q16").Value = 0.1344347627293334
$ws.Range("R16").Value = 1.209912864564
$ws.Range("S16").Value = 0.0002082189565686266
$ws.Range("T16").Value = 0.0002082189565686265
# Row 17
$ws.Range("G17").Value = 0.9167423333333334
$ws.Range("H17").Value = 2.750227
$ws.Range("I17").Value = 0.07921490167966665
$ws.Range("J17").Value = 0.07921490167966663
$ws.Range("M17").Value = 15.02284966666667
$ws.Range("N17").Value = 45.068549
$ws.Range("O17").Value = 0.2692783275177917
$ws.Range("P17").Value = 0.2692783275177917
$ws.Range("Q17").Value = 13.77208225673589
$ws.Range("R17").Value = 123.948740310623
$ws.Range("S17").Value = 0.02133085623878694
$ws.Range("T17").Value = 0.02133085623878694
# Row 18
$ws.Range("G18").Value = 0.9167423333333334
$ws.Range("H18").Value = 2.750227
$ws.Range("I18").Value = 0.07921490167966665
$ws.Range("J18").Value = 0.07921490167966663
$ws.Range("M18").Value = 5.288900666666667
$ws.Range("N18").Value = 15.866702
$ws.Range("O18").Value = 0.09480134312252211
$ws.Range("P18").Value = 0.09480134312252211
$ws.Range("Q18").Value = 4.848559137928222
$ws.Range("R18").Value = 43.63703224135401
$ws.Range("S18").Value = 0.007509679074550931
$ws.Range("T18").Value = 0.00750967907455093
# Row 19
$ws.Range("G19").Value = 0.9167423333333334
$ws.Range("H19").Value = 2.750227
$ws.Range("I19").Value = 0.07921490167966665
$ws.Range("J19").Value = 0.07921490167966663
$ws.Range("M19").Value = 0.050258
$ws.Range("N19").Value = 0.150774
$ws.Range("O19").Value = 0.0009008537317934847
$ws.Range("P19").Value = 0.0009008537317934848
$ws.Range("Q19").Value = 0.04607363618866667
$ws.Range("R19").Value = 0.414662725698
$ws.Range("S19").Value = 0.0000713610397917816826696
$ws.Range("T19").Value = 0.0000713610397917816691171
